$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Frequency"
$ws.Range("C2").Value = "monthly"
$ws.Range("C3").Value = "daily"
$ws.Range("C4").Value = "weekly"
$ws.Range("C5").Value = "daily"

$ws.Range("C5").Select()
